# Apply updated min-variance portfolio weights / returns to Sheet1.
# All touched cells are plain numeric values (no formulas in this workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Column N ("weights") updates ---
$ws.Range("N2").Value  = 0.07323693434389701
$ws.Range("N4").Value  = 0.15
$ws.Range("N5").Value  = 0.01308124912754755
$ws.Range("N6").Value  = 0.15
$ws.Range("N8").Value  = 0.02888365594205327
$ws.Range("N9").Value  = 0.1339952185821673
$ws.Range("N10").Value = 0.001
$ws.Range("N11").Value = 0.001000000000000004
$ws.Range("N13").Value = 0.07749950594290932
$ws.Range("N14").Value = 0.02112075605989106
$ws.Range("N16").Value = 0.07162963062725217
$ws.Range("N18").Value = 0.1215530493742824
$ws.Range("N19").Value = 0.15
$ws.Range("N21").Value = 0.001000000000000003

# --- Row 22 ("Portfolio return ln") updates ---
$ws.Range("B22").Value = 0.02357576643139591
$ws.Range("C22").Value = 0.01150552964432216
$ws.Range("D22").Value = 0.03118283146648781
$ws.Range("E22").Value = 0.009966260796215824
$ws.Range("F22").Value = 0.02337872376822408
$ws.Range("G22").Value = 0.04515511399120602
$ws.Range("H22").Value = -0.0275653957990564
$ws.Range("I22").Value = -0.03463897594308467
$ws.Range("J22").Value = -0.05320659953165616
$ws.Range("K22").Value = -0.004613049717445137
$ws.Range("L22").Value = -0.02494715557932752
$ws.Range("M22").Value = 0.05066715652309046

# --- Row 23 ("Portfolio return") updates ---
$ws.Range("B23").Value = 1.023855871713618
$ws.Range("C23").Value = 1.011571972827349
$ws.Range("D23").Value = 1.031674109134821
$ws.Range("E23").Value = 1.010016089370629
$ws.Range("F23").Value = 1.023654148300663
$ws.Range("G23").Value = 1.046190126049529
$ws.Range("H23").Value = 0.9728110627164643
$ws.Range("I23").Value = 0.9659540859776146
$ws.Range("J23").Value = 0.9481840978545193
$ws.Range("K23").Value = 0.9953975740541292
$ws.Range("L23").Value = 0.9753614530774264
$ws.Range("M23").Value = 1.051972692756786
$ws.Range("N23").Value = 1.125459007789521
